$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Chart")

# --- Row 4: add validationCriteria (K4) and visualisationConfig (L4) values ---
# These are new JSON config strings added to give the "Size" question
# (row 4) a visualisation config, for importer test coverage.
$ws.Range("K4").Value = '{"min": 0, "max": 300, "normalRange": {"min": 90, "max": 120}}'
$ws.Range("L4").Value = '{"yAxis":{"graphRange":{"min":40,"max":240},"interval":10}}'

# Give the new cells their own font (Calibri 12pt, black) distinct from the
# sheet's default Arial font.
$ws.Range("K4").Font.Name = "Calibri"
$ws.Range("K4").Font.Size = 12
$ws.Range("K4").Font.Color = 0

# Re-use the same computed style for L4 via a format-only copy instead of
# re-deriving the font from scratch (keeps the style table tidy).
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)

# --- Row 2: drop A2's stale explicit (but default) alignment override so it
# matches the plain font-only style used by the rest of the row (B2/C2/D2) ---
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Re-apply A2's original text (PasteSpecial(formats) only touches formatting,
# but make sure the original value is still intact).
$ws.Range("A2").Value = "PatientChartingDate"

$excel.CutCopyMode = 0
